$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1397.5385
$ws.Range("J17").Value = 1397.5385
$ws.Range("L17").Value = 4192.6155
$ws.Range("N17").Value = -4528.6155
$ws.Range("H18").Value = 1372.4445
$ws.Range("I18").Value = 1278.8572
$ws.Range("J18").Value = 1700
$ws.Range("K18").Value = 1278.8572
$ws.Range("L18").Value = 1700
$ws.Range("M18").Value = -994.8571999999999
$ws.Range("N18").Value = -2268
$ws.Range("H19").Value = 206.57143
$ws.Range("I19").Value = 211.5
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 211.5
$ws.Range("L19").Value = 200
$ws.Range("M19").Value = -36.5
$ws.Range("N19").Value = -550
$ws.Range("H70").Value = 51255.2
$ws.Range("I70").Value = 201060.4
$ws.Range("J70").Value = 1320.1333
$ws.Range("K70").Value = 603181.2
$ws.Range("L70").Value = 3960.3999
$ws.Range("M70").Value = -602911.2
$ws.Range("N70").Value = -4500.3999
$ws.Range("H73").Value = 51255.2
$ws.Range("I73").Value = 201060.4
$ws.Range("J73").Value = 1320.1333
$ws.Range("K73").Value = 603181.2
$ws.Range("L73").Value = 3960.3999
$ws.Range("M73").Value = -602245.2
$ws.Range("N73").Value = -5832.3999
$ws.Range("H94").Value = 4825.7856
$ws.Range("I94").Value = 4825.7856
$ws.Range("K94").Value = 4825.7856
$ws.Range("M94").Value = -4374.7856
$ws.Range("H99").Value = 806.3570999999999
$ws.Range("I99").Value = 637.7778
$ws.Range("J99").Value = 1109.8
$ws.Range("K99").Value = 1913.3334
$ws.Range("L99").Value = 3329.4
$ws.Range("M99").Value = -415.3334
$ws.Range("N99").Value = -6325.4
$ws.Range("H100").Value = 1850.9
$ws.Range("I100").Value = 1062.3
$ws.Range("J100").Value = 2639.5
$ws.Range("K100").Value = 1062.3
$ws.Range("L100").Value = 2639.5
$ws.Range("M100").Value = -521.3
$ws.Range("N100").Value = -3721.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3334.2173
$ws.Range("I63").Value = 3768.2307
$ws.Range("K63").Value = 3768.2307
$ws.Range("M63").Value = -3082.2307
$ws.Range("H66").Value = 3334.2173
$ws.Range("I66").Value = 3768.2307
$ws.Range("K66").Value = 18841.1535
$ws.Range("M66").Value = -15409.1535
$ws.Range("H122").Value = 1774
$ws.Range("I122").Value = 1774
$ws.Range("K122").Value = 5322
$ws.Range("M122").Value = -2872
$ws.Range("H124").Value = 23900
$ws.Range("J124").Value = 23900
$ws.Range("L124").Value = 23900
$ws.Range("N124").Value = -33720

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1689.8889
$ws.Range("I94").Value = 1218.1666
$ws.Range("K94").Value = 1218.1666
$ws.Range("M94").Value = -767.1666
$ws.Range("H99").Value = 1582.742
$ws.Range("I99").Value = 1384.9565
$ws.Range("K99").Value = 1384.9565
$ws.Range("M99").Value = 113.0435
$ws.Range("H105").Value = 4187
$ws.Range("I105").Value = 3986.25
$ws.Range("K105").Value = 3986.25
$ws.Range("M105").Value = -2239.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27030630
$ws.Range("I31").Value = 90913430
$ws.Range("J31").Value = 3292
$ws.Range("K31").Value = 90913430
$ws.Range("L31").Value = 3292
$ws.Range("M31").Value = -90913135
$ws.Range("N31").Value = -3882
$ws.Range("H34").Value = 27030630
$ws.Range("I34").Value = 90913430
$ws.Range("J34").Value = 3292
$ws.Range("K34").Value = 90913430
$ws.Range("L34").Value = 3292
$ws.Range("M34").Value = -90913228
$ws.Range("N34").Value = -3696
$ws.Range("H124").Value = 27900
$ws.Range("J124").Value = 27900
$ws.Range("L124").Value = 27900
$ws.Range("N124").Value = -32810

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1200.1818
$ws.Range("J33").Value = 1800.3334
$ws.Range("L33").Value = 10802.0004
$ws.Range("N33").Value = -11368.0004
$ws.Range("H75").Value = 4883
$ws.Range("J75").Value = 6257.143
$ws.Range("L75").Value = 18771.429
$ws.Range("N75").Value = -20767.429
$ws.Range("H78").Value = 4883
$ws.Range("J78").Value = 6257.143
$ws.Range("L78").Value = 56314.287
$ws.Range("N78").Value = -66298.287
$ws.Range("H140").Value = 2510.7222
$ws.Range("I140").Value = 989.41174
$ws.Range("J140").Value = 3871.8948
$ws.Range("K140").Value = 2968.23522
$ws.Range("L140").Value = 11615.6844
$ws.Range("M140").Value = 2211.76478
$ws.Range("N140").Value = -21975.6844

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 18055.666
$ws.Range("I97").Value = 24122.727
$ws.Range("J97").Value = 1371.25
$ws.Range("K97").Value = 24122.727
$ws.Range("L97").Value = 1371.25
$ws.Range("M97").Value = -23626.727
$ws.Range("N97").Value = -2363.25
$ws.Range("H122").Value = 2924.8696
$ws.Range("I122").Value = 3150.6667
$ws.Range("J122").Value = 2112
$ws.Range("K122").Value = 9452.000100000001
$ws.Range("L122").Value = 6336
$ws.Range("M122").Value = -7002.000100000001
$ws.Range("N122").Value = -11236
$ws.Range("H123").Value = 13036.4375
$ws.Range("J123").Value = 13036.4375
$ws.Range("L123").Value = 13036.4375
$ws.Range("N123").Value = -17936.4375

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1351.6
$ws.Range("I9").Value = 653.3333
$ws.Range("K9").Value = 653.3333
$ws.Range("M9").Value = -429.3333
$ws.Range("H40").Value = 3591
$ws.Range("I40").Value = 3986.7144
$ws.Range("J40").Value = 2667.6667
$ws.Range("K40").Value = 3986.7144
$ws.Range("L40").Value = 2667.6667
$ws.Range("M40").Value = -3850.7144
$ws.Range("N40").Value = -2939.6667
$ws.Range("H82").Value = 2283.7083
$ws.Range("I82").Value = 1828.7142
$ws.Range("K82").Value = 1828.7142
$ws.Range("M82").Value = -1467.7142
$ws.Range("H85").Value = 2283.7083
$ws.Range("I85").Value = 1828.7142
$ws.Range("K85").Value = 1828.7142
$ws.Range("M85").Value = -580.7141999999999
$ws.Range("H97").Value = 26166.166
$ws.Range("J97").Value = 26166.166
$ws.Range("L97").Value = 26166.166
$ws.Range("N97").Value = -28148.166
$ws.Range("H100").Value = 6244.4546
$ws.Range("I100").Value = 7523.625
$ws.Range("K100").Value = 7523.625
$ws.Range("M100").Value = -6982.625
$ws.Range("H101").Value = 182054
$ws.Range("J101").Value = 182054
$ws.Range("L101").Value = 182054
$ws.Range("N101").Value = -188544

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 30005.5
$ws.Range("J20").Value = 10011
$ws.Range("L20").Value = 10011
$ws.Range("N20").Value = -10491
$ws.Range("H21").Value = 30500
$ws.Range("J21").Value = 11000
$ws.Range("L21").Value = 11000
$ws.Range("N21").Value = -11470
$ws.Range("H25").Value = 29970.8
$ws.Range("J25").Value = 29970.8
$ws.Range("L25").Value = 29970.8
$ws.Range("N25").Value = -30556.8
$ws.Range("H35").Value = 30500
$ws.Range("J35").Value = 11000
$ws.Range("L35").Value = 11000
$ws.Range("N35").Value = -11580
$ws.Range("H96").Value = 1070.4286
$ws.Range("I96").Value = 942.6
$ws.Range("J96").Value = 1390
$ws.Range("K96").Value = 942.6
$ws.Range("L96").Value = 1390
$ws.Range("M96").Value = 430.4
$ws.Range("N96").Value = -4136
$ws.Range("H122").Value = 10418534
$ws.Range("I122").Value = 13159510
$ws.Range("K122").Value = 39478530
$ws.Range("M122").Value = -39476080
$ws.Range("H123").Value = 48529.42
$ws.Range("J123").Value = 48529.42
$ws.Range("L123").Value = 48529.42
$ws.Range("N123").Value = -58329.42
